$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "MISSED: $old"
    }
}

# --- Title ---
Replace-Text "Duality of Reality: Quantum Entanglement" "The Marvelous World of Biology: Unraveling the Secrets of Life"

# --- Author name: "Victor McKenzie" (1 run) -> "Dr" / "." / " Emily Carter" (3 runs) ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Text = "Dr"
$insPos = $r2.End - 1
$rDot = $d.Range($insPos, $insPos)
$rDot.InsertAfter(".")
$rDot.Font.Name = "Aptos"
$insPos2 = $insPos + 1
$rName = $d.Range($insPos2, $insPos2)
$rName.InsertAfter(" Emily Carter")

# --- Email line: 5 runs -> 3 runs ---
Replace-Text "victor" "emilycarter@edumail"
Replace-Text "mckenzie@berkeley.edu" "org"

# --- Body paragraph text replacements ---
Replace-Text "In the vast expanse of scientific exploration, the study of quantum mechanics has unveiled a captivating enigma: the phenomenon of quantum entanglement" "A journey into the realm of biology is an exploration of life itself"

Replace-Text " This perplexing phenomenon challenges our conventional notions of locality and separability, leaving us contemplating the profound interconnectedness of the universe" " From the bustling metropolis of cells within our bodies to the vibrant tapestry of ecosystems that adorn our planet, biology unveils the intricate mechanisms that govern the living world"

Replace-Text " As we delve into the mysteries of quantum entanglement, we find ourselves confronted with the fundamental question: Is reality truly separable, or is it an intricate web of interconnectedness?^l^lUnveiling the Enigmatic Dance of Entangled Particles:^l^lQuantum entanglement presents a paradoxical situation where two particles, separated by arbitrary distances, exhibit a profound correlation, defying our intuitive understanding of physical interactions" " It's a subject that ignites curiosity, challenges our understanding, and holds the key to countless mysteries that beckon us to explore further^l^lBiology unravels the enigma of life's origins, tracing the evolutionary pathways that connect all living things"

Replace-Text " This correlation manifests in the shared fate of these particles, where the measurement of one particle instantaneously influences the state of the other, regardless of the distance separating them" " It unveils the captivating dance of molecules during metabolism, revealing the symphony of reactions that sustain life"

Replace-Text " The implications of this phenomenon extend far beyond the realm of theoretical physics, inviting us to contemplate the nature of reality itself" " Through the lens of genetics, biology illuminates the blueprint of heredity, unmasking the secrets of traits passed from one generation to the next"

Replace-Text "A Deeper Dive into the Implications of Entanglement:^l^lThe study of quantum entanglement has led to profound implications for our understanding of the fundamental nature of reality" "Biology isn't just a mere collection of facts and figures; it's a dynamic and ever-evolving discipline that constantly pushes the boundaries of human knowledge"

Replace-Text " It challenges the notion of local realism, which posits that the properties of a particle are independent of measurements performed on other particles, regardless of their distance" " From the discovery of new species in the depths of rainforests to the development of groundbreaking medical treatments, biology's impact is felt in every corner of our lives"

Replace-Text " The non-local nature of entanglement suggests that the universe may be inherently interconnected, with events in one region instantaneously affecting those in another, irrespective of the vastness of the distance separating them. This concept has ignited a paradigm shift in our understanding of the universe, prompting physicists to reconsider the very fabric of reality." ""

# --- Summary paragraph ---
Replace-Text "Quantum entanglement stands as a testament to the enigmatic nature of reality" "Biology is a captivating journey into the realm of life, delving into the intricate mechanisms that govern living organisms"

Replace-Text " It presents a compelling case for the interconnectedness of the universe, defying our traditional notions of locality and separability" " From the origins of life to the evolution of species and the intricacies of genetics, biology offers a comprehensive understanding of the living world and its processes"

Replace-Text " The phenomenon of entanglement invites us to contemplate the profound depths of reality, challenging our assumptions about the fundamental nature of existence" " Its pursuit unveils the secrets of life, providing insights into human health, environmental stewardship, and the delicate balance of ecosystems"

Replace-Text " As we continue to unravel the mysteries of quantum mechanics, we may uncover even more profound insights into the fundamental nature of reality, leading us towards a deeper understanding of the universe we inhabit" " As we continue to unravel the mysteries of biology, we enrich our appreciation for the marvels of life and unlock the potential for advancements that will shape the future of our world"

# --- Append an empty paragraph at the end of the document ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

Write-Output "done"
